$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-31 Saturday" "2026-02-01 Sunday"

Replace-Text "387÷2=193, 1" "810÷5=162, 0"
Replace-Text "648÷2=324, 0" "274÷9=30, 4"
Replace-Text "456÷5=91, 1" "747÷7=106, 5"
Replace-Text "884÷3=294, 2" "186÷7=26, 4"
Replace-Text "201÷2=100, 1" "457÷5=91, 2"

Replace-Text "985÷5=197, 0" "325÷9=36, 1"
Replace-Text "440÷3=146, 2" "192÷5=38, 2"
Replace-Text "413÷7=59, 0" "462÷5=92, 2"
Replace-Text "136÷5=27, 1" "354÷7=50, 4"
Replace-Text "300÷7=42, 6" "921÷7=131, 4"

Replace-Text "425÷2=212, 1" "971÷4=242, 3"
Replace-Text "273÷6=45, 3" "235÷5=47, 0"
Replace-Text "611÷8=76, 3" "622÷2=311, 0"
Replace-Text "316÷7=45, 1" "167÷6=27, 5"
Replace-Text "532÷3=177, 1" "503÷8=62, 7"

Replace-Text "719÷4=179, 3" "176÷7=25, 1"
Replace-Text "715÷7=102, 1" "549÷4=137, 1"
Replace-Text "954÷5=190, 4" "710÷2=355, 0"
Replace-Text "758÷5=151, 3" "741÷5=148, 1"
Replace-Text "365÷8=45, 5" "432÷5=86, 2"

Replace-Text "720÷2=360, 0" "978÷3=326, 0"
Replace-Text "923÷4=230, 3" "357÷8=44, 5"
Replace-Text "350÷5=70, 0" "441÷8=55, 1"
Replace-Text "306÷8=38, 2" "612÷7=87, 3"
Replace-Text "415÷5=83, 0" "542÷6=90, 2"
